$wb = $excel.ActiveWorkbook

# --- Sheet: Trends Status ---
$wsTrends = $wb.Worksheets.Item("Trends Status")
$wsTrends.Range("A4").Value = "Insufficient Data"
$wsTrends.Range("B4").Value = 420
$wsTrends.Range("C4").Value = 300
$wsTrends.Range("A5").Value = "Trend Inconclusive"

# --- Sheet: Range Status ---
$wsRange = $wb.Worksheets.Item("Range Status")
$wsRange.Range("B2").Value = 4
$wsRange.Range("B3").Value = 47
$wsRange.Range("B4").Value = 220
$wsRange.Range("B5").Value = 363
$wsRange.Range("B6").Value = 178
$wsRange.Range("B7").Value = 131

# --- Sheet: Priority Status ---
$wsPriority = $wb.Worksheets.Item("Priority Status")
$wsPriority.Range("B2").Value = 179
$wsPriority.Range("B3").Value = 323
$wsPriority.Range("B4").Value = 441

# --- Sheet: Species qualification ---
$wsQual = $wb.Worksheets.Item("Species qualification")
$wsQual.Range("B2").Value = 943
$wsQual.Range("B5").Value = 943

# --- Sheet: SoIB-IUCN cross-tab ---
$wsCross = $wb.Worksheets.Item("SoIB-IUCN cross-tab")
$wsCross.Range("B2").Value = 14
$wsCross.Range("E2").Value = 14
$wsCross.Range("B4").Value = 42
$wsCross.Range("C4").Value = 8
$wsCross.Range("B5").Value = 17
$wsCross.Range("C5").Value = 39
$wsCross.Range("B6").Value = 91
$wsCross.Range("C6").Value = 272
$wsCross.Range("D6").Value = 423
$wsCross.Range("B7").Value = 0
$wsCross.Range("D7").Value = 4
$wsCross.Range("B8").Value = 179
$wsCross.Range("C8").Value = 441
$wsCross.Range("D8").Value = 323
$wsCross.Range("E8").Value = 943
